# Scheduled runner: refresh currentAveragePrice / LevePrice* / LeveProfit*
# columns (H:N) on each crafting-job sheet with latest market data.
# Only the price/profit columns move; leve metadata (A:G) is untouched.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 82 (ALC)
$ws.Cells.Item(82,8).Value = 2147.3333
$ws.Cells.Item(82,9).Value = 2147.3333
$ws.Cells.Item(82,11).Value = 6441.999899999999
$ws.Cells.Item(82,13).Value = -6035.999899999999

# Row 85 (ALC)
$ws.Cells.Item(85,8).Value = 2147.3333
$ws.Cells.Item(85,9).Value = 2147.3333
$ws.Cells.Item(85,11).Value = 6441.999899999999
$ws.Cells.Item(85,13).Value = -5037.999899999999

# Row 138 (ALC)
$ws.Cells.Item(138,8).Value = 2712.5
$ws.Cells.Item(138,10).Value = 3292.1738
$ws.Cells.Item(138,12).Value = 9876.5214
$ws.Cells.Item(138,14).Value = -20156.5214

$ws = $wb.Worksheets.Item("ARM")
# Row 45 (ARM)
$ws.Cells.Item(45,8).Value = 1843.75
$ws.Cells.Item(45,9).Value = 1843.75
$ws.Cells.Item(45,11).Value = 1843.75
$ws.Cells.Item(45,13).Value = -1466.75

# Row 63 (ARM)
$ws.Cells.Item(63,8).Value = 2699.5
$ws.Cells.Item(63,10).Value = 2699.5
$ws.Cells.Item(63,12).Value = 2699.5
$ws.Cells.Item(63,14).Value = -4071.5

# Row 66 (ARM)
$ws.Cells.Item(66,8).Value = 2699.5
$ws.Cells.Item(66,10).Value = 2699.5
$ws.Cells.Item(66,12).Value = 13497.5
$ws.Cells.Item(66,14).Value = -20361.5

# Row 97 (ARM)
$ws.Cells.Item(97,8).Value = 27779392
$ws.Cells.Item(97,9).Value = 41668696
$ws.Cells.Item(97,11).Value = 41668696
$ws.Cells.Item(97,13).Value = -41668200

# Row 106 (ARM)
$ws.Cells.Item(106,8).Value = 9997
$ws.Cells.Item(106,10).Value = 9997
$ws.Cells.Item(106,12).Value = 9997
$ws.Cells.Item(106,14).Value = -12521

$ws = $wb.Worksheets.Item("CRP")
# Row 4 (CRP)
$ws.Cells.Item(4,8).Value = 105
$ws.Cells.Item(4,9).Value = 105
$ws.Cells.Item(4,11).Value = 105
$ws.Cells.Item(4,13).Value = 7

# Row 99 (CRP)
$ws.Cells.Item(99,8).Value = 589439.5
$ws.Cells.Item(99,9).Value = 417856.4
$ws.Cells.Item(99,10).Value = 1001238.8
$ws.Cells.Item(99,11).Value = 417856.4
$ws.Cells.Item(99,12).Value = 1001238.8
$ws.Cells.Item(99,13).Value = -416358.4
$ws.Cells.Item(99,14).Value = -1004234.8

# Row 126 (CRP)
$ws.Cells.Item(126,8).Value = 589439.5
$ws.Cells.Item(126,9).Value = 417856.4
$ws.Cells.Item(126,10).Value = 1001238.8
$ws.Cells.Item(126,11).Value = 1253569.2
$ws.Cells.Item(126,12).Value = 3003716.4
$ws.Cells.Item(126,13).Value = -1251099.2
$ws.Cells.Item(126,14).Value = -3008656.4

$ws = $wb.Worksheets.Item("CUL")
# Row 33 (CUL)
$ws.Cells.Item(33,8).Value = 475.33334
$ws.Cells.Item(33,9).Value = 454.14285
$ws.Cells.Item(33,10).Value = 549.5
$ws.Cells.Item(33,11).Value = 2724.8571
$ws.Cells.Item(33,12).Value = 3297
$ws.Cells.Item(33,13).Value = -2441.8571
$ws.Cells.Item(33,14).Value = -3863

# Row 54 (CUL)
$ws.Cells.Item(54,8).Value = 288
$ws.Cells.Item(54,9).Value = 288
$ws.Cells.Item(54,11).Value = 864
$ws.Cells.Item(54,13).Value = -305

# Row 129 (CUL)
$ws.Cells.Item(129,8).Value = 743.875
$ws.Cells.Item(129,10).Value = 1466.5
$ws.Cells.Item(129,12).Value = 4399.5
$ws.Cells.Item(129,14).Value = -14399.5

# Row 134 (CUL)
$ws.Cells.Item(134,8).Value = 8997.5
$ws.Cells.Item(134,9).Value = 2995
$ws.Cells.Item(134,10).Value = 15000
$ws.Cells.Item(134,11).Value = 8985
$ws.Cells.Item(134,12).Value = 45000
$ws.Cells.Item(134,13).Value = -3915
$ws.Cells.Item(134,14).Value = -55140

$ws = $wb.Worksheets.Item("GSM")
# Row 5 (GSM)
$ws.Cells.Item(5,8).Value = 6004
$ws.Cells.Item(5,9).Value = 6004
$ws.Cells.Item(5,10).Value = 0
$ws.Cells.Item(5,11).Value = 6004
$ws.Cells.Item(5,12).Value = 0
$ws.Cells.Item(5,13).Value = -5892
$ws.Cells.Item(5,14).Value = $null

# Row 70 (GSM)
$ws.Cells.Item(70,8).Value = 100003200

# Row 73 (GSM)
$ws.Cells.Item(73,8).Value = 100003200

# Row 75 (GSM)
$ws.Cells.Item(75,8).Value = 53000
$ws.Cells.Item(75,10).Value = 53000
$ws.Cells.Item(75,12).Value = 53000
$ws.Cells.Item(75,14).Value = -54748

# Row 78 (GSM)
$ws.Cells.Item(78,8).Value = 53000
$ws.Cells.Item(78,10).Value = 53000
$ws.Cells.Item(78,12).Value = 159000
$ws.Cells.Item(78,14).Value = -167736

# Row 97 (GSM)
$ws.Cells.Item(97,8).Value = 976.7857
$ws.Cells.Item(97,9).Value = 976.7857
$ws.Cells.Item(97,11).Value = 976.7857
$ws.Cells.Item(97,13).Value = -480.7857

# Row 122 (GSM)
$ws.Cells.Item(122,8).Value = 5064.9473
$ws.Cells.Item(122,9).Value = 4719.647
$ws.Cells.Item(122,10).Value = 8000
$ws.Cells.Item(122,11).Value = 14158.941
$ws.Cells.Item(122,12).Value = 24000
$ws.Cells.Item(122,13).Value = -11708.941
$ws.Cells.Item(122,14).Value = -28900

$ws = $wb.Worksheets.Item("LTW")
# Row 2 (LTW)
$ws.Cells.Item(2,8).Value = 105
$ws.Cells.Item(2,10).Value = 0
$ws.Cells.Item(2,12).Value = 0
$ws.Cells.Item(2,14).Value = $null

# Row 7 (LTW)
$ws.Cells.Item(7,8).Value = 8650
$ws.Cells.Item(7,9).Value = 8650
$ws.Cells.Item(7,11).Value = 8650
$ws.Cells.Item(7,13).Value = -8538

# Row 22 (LTW)
$ws.Cells.Item(22,8).Value = 4718.6875
$ws.Cells.Item(22,9).Value = 4714.2856
$ws.Cells.Item(22,10).Value = 4722.1113
$ws.Cells.Item(22,11).Value = 4714.2856
$ws.Cells.Item(22,12).Value = 4722.1113
$ws.Cells.Item(22,13).Value = -4419.2856
$ws.Cells.Item(22,14).Value = -5312.1113

# Row 27 (LTW)
$ws.Cells.Item(27,8).Value = 4718.6875
$ws.Cells.Item(27,9).Value = 4714.2856
$ws.Cells.Item(27,10).Value = 4722.1113
$ws.Cells.Item(27,11).Value = 4714.2856
$ws.Cells.Item(27,12).Value = 4722.1113
$ws.Cells.Item(27,13).Value = -4607.2856
$ws.Cells.Item(27,14).Value = -4936.1113

# Row 55 (LTW)
$ws.Cells.Item(55,8).Value = 843.6875
$ws.Cells.Item(55,9).Value = 629.4545000000001
$ws.Cells.Item(55,10).Value = 1315
$ws.Cells.Item(55,11).Value = 629.4545000000001
$ws.Cells.Item(55,12).Value = 1315
$ws.Cells.Item(55,13).Value = -456.4545000000001
$ws.Cells.Item(55,14).Value = -1661

# Row 122 (LTW)
$ws.Cells.Item(122,8).Value = 3739.0454
$ws.Cells.Item(122,9).Value = 3077.6365
$ws.Cells.Item(122,11).Value = 9232.9095
$ws.Cells.Item(122,13).Value = -6782.9095

# Row 126 (LTW)
$ws.Cells.Item(126,8).Value = 8650
$ws.Cells.Item(126,9).Value = 8650
$ws.Cells.Item(126,11).Value = 25950
$ws.Cells.Item(126,13).Value = -23480

# Row 132 (LTW)
$ws.Cells.Item(132,8).Value = 3005
$ws.Cells.Item(132,9).Value = 0
$ws.Cells.Item(132,10).Value = 3005
$ws.Cells.Item(132,11).Value = 0
$ws.Cells.Item(132,12).Value = 9015
$ws.Cells.Item(132,13).Value = $null
$ws.Cells.Item(132,14).Value = -14075

# Row 136 (LTW)
$ws.Cells.Item(136,8).Value = 1165
$ws.Cells.Item(136,9).Value = 325
$ws.Cells.Item(136,10).Value = 2005
$ws.Cells.Item(136,11).Value = 975
$ws.Cells.Item(136,12).Value = 6015
$ws.Cells.Item(136,13).Value = 1575
$ws.Cells.Item(136,14).Value = -11115

$ws = $wb.Worksheets.Item("WVR")
# Row 2 (WVR)
$ws.Cells.Item(2,8).Value = 3253.9167
$ws.Cells.Item(2,9).Value = 105
$ws.Cells.Item(2,10).Value = 18998.5
$ws.Cells.Item(2,11).Value = 105
$ws.Cells.Item(2,12).Value = 18998.5
$ws.Cells.Item(2,13).Value = 7
$ws.Cells.Item(2,14).Value = -19222.5

# Row 13 (WVR)
$ws.Cells.Item(13,8).Value = 0
$ws.Cells.Item(13,9).Value = 0
$ws.Cells.Item(13,11).Value = 0
$ws.Cells.Item(13,13).Value = $null

# Row 81 (WVR)
$ws.Cells.Item(81,8).Value = 2487.5557
$ws.Cells.Item(81,9).Value = 2148.1667
$ws.Cells.Item(81,10).Value = 3166.3333
$ws.Cells.Item(81,11).Value = 4296.3334
$ws.Cells.Item(81,12).Value = 6332.6666
$ws.Cells.Item(81,13).Value = -3235.3334
$ws.Cells.Item(81,14).Value = -8454.6666

# Row 84 (WVR)
$ws.Cells.Item(84,8).Value = 2487.5557
$ws.Cells.Item(84,9).Value = 2148.1667
$ws.Cells.Item(84,10).Value = 3166.3333
$ws.Cells.Item(84,11).Value = 21481.667
$ws.Cells.Item(84,12).Value = 31663.333
$ws.Cells.Item(84,13).Value = -16177.667
$ws.Cells.Item(84,14).Value = -42271.333

# Row 126 (WVR)
$ws.Cells.Item(126,8).Value = 6176.75
$ws.Cells.Item(126,9).Value = 4900.6665
$ws.Cells.Item(126,11).Value = 14701.9995
$ws.Cells.Item(126,13).Value = -12231.9995

